$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at F (existing F,G shift to H,I, inheriting their
# formats/values/formulas automatically).
$ws.Columns("F:G").Insert()

# The new F/G columns should carry the "standard body cell" style used by
# most of the table (same as column A/E use, style index 8 / border id 2).
# Grab that formatting from an existing body cell (A3) and stamp it across
# the new columns first; row 2 keeps the slightly different header-adjacent
# style it already inherited from the insert (matches target diff).
$fmtSrc = $ws.Range("A3")
$fmtSrc.Copy()
$ws.Range("F3:G37").PasteSpecial(-4122)

# Header labels for the two new columns.
$ws.Range("G1").Value = "지역"
$ws.Range("F1").Value = "준공년도"

# Data values - every site row gets the same region / completion year.
$ws.Range("G2:G37").Value = "수도권"
$ws.Range("F2:F37").Value = "2024년"

# Column widths approximating the target layout (9.25 / 5.5 characters).
$ws.Columns("F").ColumnWidth = 8.571428571428571
$ws.Columns("G").ColumnWidth = 4.857142857142857

# Sheet view: zoom to 85%, select F3:F37 (mirrors the authored edit).
$excel.ActiveWindow.Zoom = 85
$ws.Range("F3:F37").Select()
